$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q1" sheet right before "总计", and fill it
#    with the per-fund holding breakdown (same layout as the other
#    quarterly sheets, e.g. "2021-Q4").
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Worksheet handles are position-based, so re-resolve "总计" by name now
# that the new sheet has shifted its index.
$totalSheet = $wb.Worksheets.Item("总计")

# Reuse the header row's formatting (bold, centered, bordered) from the
# template sheet so the new sheet's style indices line up with the rest
# of the workbook instead of minting new ones.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# B:G hold text-like values (fund codes with leading zeros, decimal
# numbers formatted as strings) in the source data, so force text
# format before assigning - otherwise Excel coerces them to numbers
# and leading zeros are lost.
$newSheet.Range("B1:G12").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "515900", "博时中证央企创新驱动ETF",                          "47.96", "99.42", "2.89",  "1.3860", 10),
    @(1,  "516950", "银华中证基建交易型开放式指数证券投资基金",          "10.41", "97.55", "8.41",  "0.8755", 1),
    @(2,  "515680", "嘉实中证央企创新驱动ETF",                          "17.86", "99.22", "2.90",  "0.5179", 9),
    @(3,  "515600", "广发中证央企创新驱动ETF",                          "17.26", "99.02", "2.89",  "0.4988", 9),
    @(4,  "159974", "富国中证央企创新驱动ETF",                          "5.37",  "99.51", "2.89",  "0.1552", 9),
    @(5,  "160639", "鹏华中证高铁产业指数（LOF）",                       "0.89",  "94.72", "12.86", "0.1145", 3),
    @(6,  "160638", "鹏华中证一带一路主题指数（LOF）",                   "3.52",  "94.71", "2.75",  "0.0968", 9),
    @(7,  "011243", "万家惠裕回报6个月持有期混合型证券投资基金A",       "4.93",  "23.04", "1.21",  "0.0597", 3),
    @(8,  "167503", "安信中证一带一路主题指数",                         "1.09",  "94.39", "2.74",  "0.0299", 9),
    @(9,  "006478", "长盛多因子策略优选股票",                           "0.51",  "84.41", "4.69",  "0.0239", 3),
    @(10, "011244", "万家惠裕回报6个月持有期混合型证券投资基金C",       "0.14",  "23.04", "1.21",  "0.0017", 3)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing history down by one row, and renumber the index column.
# ---------------------------------------------------------------------
$totalSheet.Range("A2").EntireRow.Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 3.76

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
